$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A5").Value = 9.24
$wsSummary.Range("E5").Value = 9.24

# ---------------------------------------------------------------------------
# Sheet: Original Schedule
# ---------------------------------------------------------------------------
$wsOrig = $wb.Worksheets.Item("Original Schedule")
$wsOrig.Range("B5").Value = 853.95
$wsOrig.Range("C5").Value = 2459.94
$wsOrig.Range("D5").Value = 33.770000000000003
$wsOrig.Range("F5").Value = 6.08
$wsOrig.Range("G5").Value = 893.8
$wsOrig.Range("C6").Value = 1596.47
$wsOrig.Range("C7").Value = 725.01
$wsOrig.Range("B8").Value = 725.01
$wsOrig.Range("G8").Value = 732.15

# ---------------------------------------------------------------------------
# Sheet: Repayment schedule
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Range("J5").Value = 6.08
$wsRepay.Range("K5").Value = 893.8
$wsRepay.Range("Q5").Value = 893.8

# ---------------------------------------------------------------------------
# Sheet: Transactions
# ---------------------------------------------------------------------------
$wsTxn = $wb.Worksheets.Item("Transactions")
$wsTxn.Range("A2").Value = 466
$wsTxn.Range("E2").Value = 58
$wsTxn.Range("J2").Value = 5151.97
$wsTxn.Range("A3").Value = 465
$wsTxn.Range("E3").Value = 6.08
$wsTxn.Range("A4").Value = 457
$wsTxn.Range("A5").Value = 449
$wsTxn.Range("A6").Value = 441
$wsTxn.Range("A7").Value = 462
$wsTxn.Range("A8").Value = 461
$wsTxn.Range("A9").Value = 455
$wsTxn.Range("A10").Value = 447
$wsTxn.Range("A11").Value = 439
$wsTxn.Range("A12").Value = 438
$wsTxn.Range("A13").Value = 437
$wsTxn.Range("A14").Value = 436

# ---------------------------------------------------------------------------
# Sheet: ChargesTab
# ---------------------------------------------------------------------------
$wsCharges = $wb.Worksheets.Item("ChargesTab")

# Update the overdue-fee charge amount text (shared by G6 and J6).
$wsCharges.Range("G6").Value = "$6.08"
$wsCharges.Range("J6").Value = "$6.08"

# Re-apply the plain (non-bold-lookalike) font formatting used elsewhere in
# the workbook to the data rows, picking it up from cells that already carry
# the desired style so no new styles are introduced.
$wsRepay.Range("A2").Copy()
$wsCharges.Range("A3:D6").PasteSpecial(-4122)
$wsCharges.Range("F3:J6").PasteSpecial(-4122)

$wsRepay.Range("C2").Copy()
$wsCharges.Range("E3:E6").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Selections (also drives which sheet tab / cell is active when saved)
# ---------------------------------------------------------------------------
$wsSummary.Range("D10").Select()
$wsOrig.Range("F13").Select()
$wsRepay.Range("L5").Select()
$wsTxn.Range("J5").Select()
$wsCharges.Range("M7").Select()
